$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 804, pushing existing rows 804-866 down to 805-867.
$ws.Rows.Item(804).Insert()

# Populate the newly inserted row 804 with the new record.
# Columns A,B,C,E,F,G,H,I,O,R keep the same values as the (now shifted) row 805,
# i.e. the same as the original row 804 had before the insert.
$ws.Cells.Item(804, 1).Value = 3
$ws.Cells.Item(804, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(804, 3).Value = "Coquimbo"
$ws.Cells.Item(804, 4).Value = 45223
$ws.Cells.Item(804, 5).Value = 5
$ws.Cells.Item(804, 6).Value = 100112032
$ws.Cells.Item(804, 7).Value = "Zapallo italiano"
$ws.Cells.Item(804, 8).Value = "Sin especificar"
$ws.Cells.Item(804, 9).Value = "Primera"
$ws.Cells.Item(804, 10).Value = 50
$ws.Cells.Item(804, 11).Value = 15000
$ws.Cells.Item(804, 12).Value = 15000
$ws.Cells.Item(804, 13).Value = 15000
$ws.Cells.Item(804, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(804, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(804, 16).Value = 250
$ws.Cells.Item(804, 17).Value = 60
$ws.Cells.Item(804, 18).Value = "Hortaliza"
